$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Tank Diameter -> Viscosity (Step column F2 stays blank, untouched)
$ws.Cells.Item(2, 1).Value = "Viscosity"
$ws.Cells.Item(2, 2).Value = 800
$ws.Cells.Item(2, 3).Value = 50
$ws.Cells.Item(2, 4).Value = 882.2426813475736
$ws.Cells.Item(2, 5).Value = 717.7573186524264
$ws.Cells.Item(2, 7).Value = 0.9
$ws.Cells.Item(2, 8).Value = "Continuous"

# Row 3: Tank Level -> Densidade (Step column F3 stays blank, untouched)
$ws.Cells.Item(3, 1).Value = "Densidade"
$ws.Cells.Item(3, 2).Value = 1500
$ws.Cells.Item(3, 3).Value = 194.1122415647322
$ws.Cells.Item(3, 4).Value = 2000
$ws.Cells.Item(3, 5).Value = 1000
$ws.Cells.Item(3, 7).Value = 0.99
$ws.Cells.Item(3, 8).Value = "Continuous"

# Row 4: Viscosity -> Numero de Pratos (now Discrete with Step = 1)
$ws.Cells.Item(4, 1).Value = "Numero de Pratos"
$ws.Cells.Item(4, 2).Value = 0
$ws.Cells.Item(4, 3).Value = 1
$ws.Cells.Item(4, 4).Value = 15
$ws.Cells.Item(4, 5).Value = 10
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.95
$ws.Cells.Item(4, 8).Value = "Discrete"

# Row 5: Number of Impellers -> Discreto 2 (Discrete with Step = 0.5)
$ws.Cells.Item(5, 1).Value = "Discreto 2"
$ws.Cells.Item(5, 2).Value = 0
$ws.Cells.Item(5, 3).Value = 1
$ws.Cells.Item(5, 4).Value = 5
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 0.5
$ws.Cells.Item(5, 7).Value = 0.95
$ws.Cells.Item(5, 8).Value = "Discrete"

# Remove old rows 6-9 (Impeller 1/2/3 Angle, Densidade) - no longer used
$ws.Range("A6:H9").Clear()
